# Updates the cryptocurrency price/volume table (columns D = Price, E = Volume(1h))
# to the latest scraped values, per the "Updated cryptos list ... GitHub Actions" run.
#
# NOTE: some Price values are plain numeric-looking strings (e.g. "97.09").
# Assigning such a string straight to .Value would make Excel auto-convert the
# cell to a real Number, which would both change the cell's stored type and can
# silently drop meaningful trailing zeros (e.g. "0.0700" -> 0.07). To keep these
# cells as text (matching their original t="inlineStr"/text representation) we
# prefix the value with a leading apostrophe, exactly like typing it by hand in
# Excel, and then reset the cell Style back to "Normal" so no stray number
# format/style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Leading apostrophe forces Excel to store the value as text even when it
    # looks like a number.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.936.24"
$ws.Range("E2").Value = "  +0.74%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.298.61"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "299.26"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "97.09"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.54%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.20%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "33.75"
$ws.Range("E10").Value = "  +0.89%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.46%  "

# Row 12 - OKB
Set-TextValue $ws.Range("D12") "48.76"
$ws.Range("E12").Value = "  -2.98%  "

# Row 13 - TRON
Set-TextValue $ws.Range("D13") "0.115"
$ws.Range("E13").Value = "  +2.81%  "

# Row 14 - Chainlink
$ws.Range("E14").Value = "  +12.35%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +1.60%  "

# Row 16 - Wrapped liquid staked Ether 2.0
$ws.Range("D16").Value = "2.650.20"
$ws.Range("E16").Value = "  +0.48%  "

# Row 17 - Wrapped Ether
$ws.Range("D17").Value = "2.295.16"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +2.55%  "

# Row 19 - Wrapped BTC
$ws.Range("D19").Value = "42.870.91"
$ws.Range("E19").Value = "  +0.85%  "

# Row 20 - Internet Computer (DFINITY)
Set-TextValue $ws.Range("D20") "11.62"
$ws.Range("E20").Value = "  +1.31%  "

# Row 21 - Shiba Inu
$ws.Range("E21").Value = "  +0.68%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.55%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "67.44"
$ws.Range("E23").Value = "  +1.03%  "

# Row 24 - Bitcoin Cash
Set-TextValue $ws.Range("D24") "236.57"
$ws.Range("E24").Value = "  +0.78%  "

# Row 25 - Immutable X
Set-TextValue $ws.Range("D25") "2.03"
$ws.Range("E25").Value = "  +5.00%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  -1.58%  "

# Row 28 - Ethereum Classic
Set-TextValue $ws.Range("D28") "24.33"
$ws.Range("E28").Value = "  -0.61%  "

# Row 29 - Monero
Set-TextValue $ws.Range("D29") "166.78"
$ws.Range("E29").Value = "  +0.34%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +0.26%  "

# Row 31 - Injective Protocol
Set-TextValue $ws.Range("D31") "33.76"
$ws.Range("E31").Value = "  +0.09%  "

# Row 32 - Cosmos
Set-TextValue $ws.Range("D32") "9.10"

# Row 33 - First Digital USD
$ws.Range("E33").Value = "  +0.08%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  -0.20%  "

# Row 35 - Render Token
$ws.Range("E35").Value = "  +5.81%  "

# Row 36 - WEMIX Token
$ws.Range("E36").Value = "  +1.76%  "

# Row 37 - Celestia
Set-TextValue $ws.Range("D37") "16.78"
$ws.Range("E37").Value = "  +3.36%  "

# Row 38 - Hedera
Set-TextValue $ws.Range("D38") "0.0700"
$ws.Range("E38").Value = "  +0.48%  "

# Row 39 - Lido DAO Token
$ws.Range("E39").Value = "  -0.11%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +0.53%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +0.15%  "

# Row 43 - ApeX Protocol
$ws.Range("E43").Value = "  -4.40%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.991.32"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +0.59%  "

# Row 46 - FraxShare
Set-TextValue $ws.Range("D46") "9.83"
$ws.Range("E46").Value = "  +1.41%  "

# Row 47 - EnergySwap
Set-TextValue $ws.Range("D47") "17.50"
$ws.Range("E47").Value = "  -1.95%  "

# Row 48 - NEAR Protocol
$ws.Range("E48").Value = "  +0.14%  "

# Row 49 - Rocket Pool ETH
$ws.Range("D49").Value = "2.528.86"
$ws.Range("E49").Value = "  +0.92%  "

# Row 50 - MultiversX
Set-TextValue $ws.Range("D50") "53.16"
$ws.Range("E50").Value = "  -0.01%  "

# Row 51 - THORChain
$ws.Range("E51").Value = "  -1.61%  "
